# The deck's main (slide-master) theme was switched from the custom
# "Integral" theme to the stock default "Office Theme" palette, while the
# previous "Integral" colours were pushed into the secondary theme part
# (the one now used only by the Notes Master).
#
# We reach the live theme's 12-slot DrawingML colour scheme through
# Slide.ThemeColorScheme (Slide -> Master -> Theme -> clrScheme), and
# overwrite every slot with the standard Office theme RGB values. RGB is
# expressed the VBA way: a single long built from B*65536 + G*256 + R.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# index -> (scheme slot, target colour as VBA RGB long == B*65536+G*256+R)
$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72
